$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.600.93"
$ws.Range("E2").Value = "'  -0.28%  "

$ws.Range("D3").Value = "'1.842.04"
$ws.Range("E3").Value = "'  -0.39%  "

$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "'  -2.31%  "

$ws.Range("D5").Value = "'316.80"
$ws.Range("E5").Value = "'  -1.66%  "

$ws.Range("E6").Value = "'  -2.19%  "

$ws.Range("D7").Value = "'0.4299"
$ws.Range("E7").Value = "'  -1.95%  "

$ws.Range("D8").Value = "'0.3732"
$ws.Range("E8").Value = "'  -1.54%  "

$ws.Range("D9").Value = "'0.07293"
$ws.Range("E9").Value = "'  -1.17%  "

$ws.Range("D10").Value = "'0.8709"
$ws.Range("E10").Value = "'  -1.16%  "

$ws.Range("D11").Value = "'21.38"
$ws.Range("E11").Value = "'  -0.66%  "

$ws.Range("D12").Value = "'1.846.32"
$ws.Range("E12").Value = "'  -0.32%  "

$ws.Range("D13").Value = "'6.710"
$ws.Range("E13").Value = "'  +0.11%  "

$ws.Range("D14").Value = "'5.392"
$ws.Range("E14").Value = "'  -1.92%  "

$ws.Range("D15").Value = "'0.07110"
$ws.Range("E15").Value = "'  -0.46%  "

$ws.Range("D16").Value = "'88.54"
$ws.Range("E16").Value = "'  +4.24%  "

$ws.Range("E17").Value = "'  -2.38%  "

$ws.Range("D18").Value = "'0.000008969"
$ws.Range("E18").Value = "'  -0.90%  "

$ws.Range("E19").Value = "'  -2.07%  "

$ws.Range("D20").Value = "'15.33"
$ws.Range("E20").Value = "'  -0.70%  "

$ws.Range("D21").Value = "'27.606.62"
$ws.Range("E21").Value = "'  -0.28%  "

$ws.Range("D22").Value = "'5.189"
$ws.Range("E22").Value = "'  -1.79%  "

$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = "'  -2.55%  "

$ws.Range("D24").Value = "'2.071.41"
$ws.Range("E24").Value = "'  -0.60%  "

$ws.Range("D25").Value = "'1.967"
$ws.Range("E25").Value = "'  -3.87%  "

$ws.Range("D26").Value = "'154.54"
$ws.Range("E26").Value = "'  -2.18%  "

$ws.Range("D27").Value = "'18.46"
$ws.Range("E27").Value = "'  -0.98%  "

$ws.Range("D28").Value = "'2.151"
$ws.Range("E28").Value = "'  +7.66%  "

$ws.Range("D29").Value = "'5.314"
$ws.Range("E29").Value = "'  +0.01%  "

$ws.Range("D30").Value = "'117.49"
$ws.Range("E30").Value = "'  -0.15%  "

$ws.Range("D31").Value = "'0.08903"
$ws.Range("E31").Value = "'  -1.48%  "

$ws.Range("D32").Value = "'1.213"
$ws.Range("E32").Value = "'  +0.67%  "

$ws.Range("D33").Value = "'0.7720"
$ws.Range("E33").Value = "'  +0.40%  "

$ws.Range("D34").Value = "'4.516"
$ws.Range("E34").Value = "'  -0.64%  "

$ws.Range("D35").Value = "'2.895"
$ws.Range("E35").Value = "'  -3.42%  "

$ws.Range("D36").Value = "'1.007"
$ws.Range("E36").Value = "'  -2.29%  "

$ws.Range("D37").Value = "'1.126"
$ws.Range("E37").Value = "'  -2.07%  "

$ws.Range("E38").Value = "'  -0.23%  "

$ws.Range("D39").Value = "'0.05293"
$ws.Range("E39").Value = "'  +0.67%  "

$ws.Range("B40").Value = "'MXToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.894"
$ws.Range("E40").Value = "'  +1.80%  "

$ws.Range("B41").Value = "'FraxShare"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.139"
$ws.Range("E41").Value = "'  +4.16%  "

$ws.Range("D42").Value = "'0.1686"
$ws.Range("E42").Value = "'  +1.04%  "

$ws.Range("D43").Value = "'0.5107"
$ws.Range("E43").Value = "'  -1.22%  "

$ws.Range("D44").Value = "'8.741"
$ws.Range("E44").Value = "'  +0.17%  "

$ws.Range("D45").Value = "'10.70"
$ws.Range("E45").Value = "'  -0.10%  "

$ws.Range("D46").Value = "'106.82"
$ws.Range("E46").Value = "'  -3.03%  "

$ws.Range("D47").Value = "'0.4738"
$ws.Range("E47").Value = "'  +1.01%  "

$ws.Range("D48").Value = "'0.06446"
$ws.Range("E48").Value = "'  -2.72%  "

$ws.Range("D49").Value = "'1.007"
$ws.Range("E49").Value = "'  -2.38%  "

$ws.Range("D50").Value = "'1.680"
$ws.Range("E50").Value = "'  -1.06%  "

$ws.Range("D51").Value = "'1.841"
$ws.Range("E51").Value = "'  -2.73%  "
